$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BD")

# --- Row 1 extra annotation cells ---------------------------------------
# Order of assignment controls shared-string index allocation, so these
# are written in the exact sequence the original author must have used.
$ws.Range("B1").Value = "high var means DIC =[0.1, 1.0, 1.9] of nomial value"
$ws.Range("F1").Value = "low var:DIC=[0.3, 1.0, 1.3] "
$ws.Range("A1").Value = "5 years, 5 days low ng price"

# --- Row 10 extra annotation cells --------------------------------------
$ws.Range("A10").Value = "more tests with medim carbon tax "
$ws.Range("H1").Value = "coal-new-igcc"
$ws.Range("C10").Value = "coal-igcc-css-new"

# --- Row 11 header (mirrors row 2's LP-table header, plus a new H col) --
$ws.Range("A11").Value = "instance"
$ws.Range("B11").Value = "LB"
$ws.Range("C11").Value = "UB"
$ws.Range("D11").Value = "Gap"
$ws.Range("E11").Value = "Benders time"
$ws.Range("F11").Value = "Heuristic time"
$ws.Range("G11").Value = "Wall time"

# --- Row 12 ---------------------------------------------------------------
$ws.Range("A12").Value = "DIC= [0.2 1.0, 1.3], high ng"
$ws.Range("B12").Formula = "=77997.32707/1000"
$ws.Range("C12").Formula = "=78949.46069/1000"
$ws.Range("D12").Formula = "=(C12-B12)/C12"
$ws.Range("E12").Value = 43200
$ws.Range("F12").Value = 273
$ws.Range("G12").Value = 44970

# --- Row 13 -----------------------------------------------------------------
$ws.Range("A13").Value = "DIC= [0.1 1.0, 1.3], high ng"

$ws.Range("H11").Value = "solution description"

$ws.Range("B13").Formula = "=76629/1000"
$ws.Range("C13").Value = 77.6
$ws.Range("D13").Formula = "=(C13-B13)/C13"
$ws.Range("E13").Value = 43200
$ws.Range("F13").Value = 291
$ws.Range("G13").Value = 45036

# --- Row 14 -----------------------------------------------------------------
$ws.Range("A14").Value = "DIC= [0.2 1.0, 1.3], high ng * 1.5"
$ws.Range("B14").Formula = "=86701/1000"
$ws.Range("C14").Value = 87.335
$ws.Range("D14").Formula = "=(C14-B14)/C14"
$ws.Range("E14").Value = 43200
$ws.Range("F14").Value = 282
$ws.Range("G14").Value = 45047

# --- Row 15 -----------------------------------------------------------------
$ws.Range("A15").Value = "DIC= [0.5 1.0, 1.3], high ng * 1.5"
$ws.Range("B15").Value = 90.274
$ws.Range("C15").Value = 90.842
$ws.Range("D15").Formula = "=(C15-B15)/C15"
$ws.Range("E15").Value = 43200
$ws.Range("F15").Value = 89
$ws.Range("G15").Value = 44763

$ws.Range("H12").Value = "1 coal plant installed at t=1. Continue installation for low price scenario"
$ws.Range("H13").Value = "1 coal plant installed at t=1. Continue installation for low price scenario"
$ws.Range("H14").Value = "1 coal plant installed at t=1. Continue installation for low price scenario"
$ws.Range("E10").Value = "time limit: 43200 secs"
$ws.Range("H15").Value = "no coal plant installed in all scenarios"

Write-Host "step1 ok"
